# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型"
# sheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (F column: row -> new value)
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 4705
    3  = 1867
    6  = 3160
    9  = 281
    10 = 646
    13 = 405
    14 = 140
    15 = 1795
    16 = 1379
    18 = 1642
    19 = 23
    20 = 130
    21 = 618
    32 = 3986
    34 = 781
    36 = 1463
    38 = 1886
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# Sheet "全部类型" (F column: row -> new value)
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 4705
    3  = 1867
    6  = 3160
    9  = 281
    10 = 646
    14 = 405
    15 = 140
    16 = 1795
    17 = 1379
    19 = 1642
    20 = 23
    21 = 130
    22 = 618
    33 = 3986
    37 = 781
    39 = 1463
    41 = 1886
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
